# New website version v3
#
# Updates the "professors_promotion_type" data sheet:
#   - Row 21 is corrected from year 2001/century 21 to year 1963/century 20.
#   - Rows 22-56 are newly appended with additional promotion-type records,
#     extending the table through year 2012.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# promotion type, count, year (kept as text, matching the existing "year"
# column which stores 4-digit years as text rather than numbers), century
$data = @(
    @("Promotie", 1,  "1963", 20),
    @("Promotie", 1,  "1977", 20),
    @("Promotie", 2,  "1978", 20),
    @("Promotie", 5,  "1980", 20),
    @("Promotie", 2,  "1981", 20),
    @("Promotie", 5,  "1982", 20),
    @("Promotie", 4,  "1983", 20),
    @("Promotie", 10, "1984", 20),
    @("Promotie", 13, "1985", 20),
    @("Promotie", 8,  "1986", 20),
    @("Promotie", 16, "1987", 20),
    @("Promotie", 17, "1988", 20),
    @("Promotie", 20, "1989", 20),
    @("Promotie", 11, "1990", 20),
    @("Promotie", 15, "1991", 20),
    @("Promotie", 20, "1992", 20),
    @("Promotie", 11, "1993", 20),
    @("Promotie", 17, "1994", 20),
    @("Promotoe", 1,  "1994", 20),
    @("Promotie", 9,  "1995", 20),
    @("Promotie", 7,  "1996", 20),
    @("Promotie", 11, "1997", 20),
    @("Promotie", 5,  "1998", 20),
    @("Promotie", 4,  "1999", 20),
    @("Promotie", 9,  "2000", 21),
    @("Promotie", 1,  "2001", 21),
    @("Promotie", 4,  "2002", 21),
    @("Promotie", 5,  "2003", 21),
    @("Promotie", 3,  "2004", 21),
    @("Promotie", 1,  "2005", 21),
    @("Promotie", 1,  "2006", 21),
    @("Promotie", 3,  "2007", 21),
    @("Promotie", 3,  "2008", 21),
    @("Promotie", 2,  "2010", 21),
    @("Promotie", 2,  "2011", 21),
    @("Promotie", 1,  "2012", 21)
)

$startRow = 21
$endRow = $startRow + $data.Length - 1

# Keep the "year" column (C) formatted/stored as text for every row we
# touch, otherwise 4-digit year strings like "1963" would be auto-detected
# as numbers when assigned.
$ws.Range("C$startRow`:C$endRow").NumberFormat = "@"

$row = $startRow
foreach ($rec in $data) {
    $ws.Cells.Item($row, 1).Value = $rec[0]
    $ws.Cells.Item($row, 2).Value = $rec[1]
    $ws.Cells.Item($row, 3).Value = $rec[2]
    $ws.Cells.Item($row, 4).Value = $rec[3]
    $row++
}
